$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 170, pushing existing rows 170-228 down to 172-230.
$ws.Rows.Item(170).Insert()
$ws.Rows.Item(171).Insert()

# Fill in the new row 170 with the new weekly price entry.
$ws.Range("A170").Value = 4
$ws.Range("B170").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C170").Value = "Los Lagos"
$ws.Range("D170").Value = 44663
$ws.Range("E170").Value = 10
$ws.Range("F170").Value = 100112032
$ws.Range("G170").Value = "Zapallo italiano"
$ws.Range("H170").Value = "Sin especificar"
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 120
$ws.Range("K170").Value = 15000
$ws.Range("L170").Value = 15000
$ws.Range("M170").Value = 15000
$ws.Range("N170").Value = "$/caja 50 unidades"
$ws.Range("O170").Value = "Región Metropolitana"
$ws.Range("P170").Value = 300
$ws.Range("Q170").Value = 50
$ws.Range("R170").Value = "Hortaliza"

# Fill in the new row 171 with the second new weekly price entry.
$ws.Range("A171").Value = 4
$ws.Range("B171").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C171").Value = "Los Lagos"
$ws.Range("D171").Value = 44663
$ws.Range("E171").Value = 10
$ws.Range("F171").Value = 100112032
$ws.Range("G171").Value = "Zapallo italiano"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 120
$ws.Range("K171").Value = 13000
$ws.Range("L171").Value = 13000
$ws.Range("M171").Value = 13000
$ws.Range("N171").Value = "$/caja 50 unidades"
$ws.Range("O171").Value = "Región del Maule"
$ws.Range("P171").Value = 260
$ws.Range("Q171").Value = 50
$ws.Range("R171").Value = "Hortaliza"
